# Euclides.xlsx - "Listas sem duplicação de professores"
# Remove duplicated teacher entries across the weekly schedule grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value  = "[-, 'MCT-3A-Tecnologia da soldagem', -, -]"

$ws.Range("B3").Value  = "-"
$ws.Range("E3").Value  = "[-, 'MCT-3A-Tecnologia da soldagem', -, -]"

$ws.Range("B4").Value  = "-"
$ws.Range("D4").Value  = "-"

$ws.Range("B6").Value  = "[-, -, 'MCT-3A-Tecnologia da soldagem', -]"
$ws.Range("D6").Value  = "-"

$ws.Range("B7").Value  = "[-, -, 'MCT-3A-Tecnologia da soldagem', -]"

$ws.Range("E8").Value  = "-"

$ws.Range("C18").Value = "-"

$ws.Range("C19").Value = "-"

$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"

$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "-"
